# Update symbol list values (cryptos.xlsx) - GitHub Actions scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new value (kept as text to preserve exact formatting,
# e.g. trailing zeros such as "243.43" or "0.8711")
$updates = [ordered]@{
    "D2"  = "243.43"
    "D4"  = "5.304"
    "D5"  = "0.05795"
    "D6"  = "6.494"
    "D7"  = "3.335"
    "D9"  = "0.8711"
    "D10" = "0.1384"
    "D11" = "0.07279"
    "D12" = "0.03075"
    "D13" = "0.03059"
    "D14" = "0.09313"
    "D15" = "3.851"
    "D16" = "0.001537"
    "D17" = "0.04713"
    "D18" = "0.0006042"
    "E18" = "17OneONE"
    "D19" = "0.006056"
    "D21" = "0.004591"
    "D22" = "0.00008702"
    "D24" = "2.147"
    "D28" = "0.0002344"
    "D40" = "0.03782"
    "D41" = "0.006352"
    "D42" = "0.1052"
    "D43" = "0.002701"
    "D44" = "0.007003"
    "D45" = "0.00005489"
    "D47" = "0.5501"
    "E47" = "46CoinbaseStockTokenCOINWorstin24h"
    "D48" = "0.006579"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
